# Apply the updated NATMI TPM recomputation values to rows 2-10 of Sheet1.
# (Sending cluster / Ligand / Receptor / Target cluster columns A-D are unchanged.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04027
$ws.Range("H2").Value = 0.12081
$ws.Range("I2").Value = 0.01318991723029425
$ws.Range("J2").Value = 0.01318991723029425
$ws.Range("M2").Value = 72.266001
$ws.Range("N2").Value = 216.798003
$ws.Range("O2").Value = 0.2949652269937106
$ws.Range("P2").Value = 0.2949652269937106
$ws.Range("Q2").Value = 2.91015186027
$ws.Range("R2").Value = 26.19136674243
$ws.Range("S2").Value = 0.003890566929861999
$ws.Range("T2").Value = 0.003890566929861998

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04027
$ws.Range("H3").Value = 0.12081
$ws.Range("I3").Value = 0.01318991723029425
$ws.Range("J3").Value = 0.01318991723029425
$ws.Range("N3").Value = 410.023338
$ws.Range("O3").Value = 0.5578585839920717
$ws.Range("P3").Value = 0.5578585839920718
$ws.Range("Q3").Value = 5.503879940419999
$ws.Range("R3").Value = 49.53491946378
$ws.Range("S3").Value = 0.00735810854906458
$ws.Range("T3").Value = 0.007358108549064581

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04027
$ws.Range("H4").Value = 0.12081
$ws.Range("I4").Value = 0.01318991723029425
$ws.Range("J4").Value = 0.01318991723029425
$ws.Range("M4").Value = 36.057927
$ws.Range("N4").Value = 108.173781
$ws.Range("O4").Value = 0.1471761890142177
$ws.Range("P4").Value = 0.1471761890142177
$ws.Range("Q4").Value = 1.45205272029
$ws.Range("R4").Value = 13.06847448261
$ws.Range("S4").Value = 0.001941241751367674
$ws.Range("T4").Value = 0.001941241751367673

# Row 5
$ws.Range("G5").Value = 0.9943730000000001
$ws.Range("I5").Value = 0.3256940046198011
$ws.Range("J5").Value = 0.325694004619801
$ws.Range("M5").Value = 72.266001
$ws.Range("N5").Value = 216.798003
$ws.Range("O5").Value = 0.2949652269937106
$ws.Range("P5").Value = 0.2949652269937106
$ws.Range("Q5").Value = 71.859360212373
$ws.Range("R5").Value = 646.734241911357
$ws.Range("S5").Value = 0.09606840600317025
$ws.Range("T5").Value = 0.09606840600317024

# Row 6
$ws.Range("G6").Value = 0.9943730000000001
$ws.Range("I6").Value = 0.3256940046198011
$ws.Range("J6").Value = 0.325694004619801
$ws.Range("N6").Value = 410.023338
$ws.Range("O6").Value = 0.5578585839920717
$ws.Range("P6").Value = 0.5578585839920718
$ws.Range("S6").Value = 0.1816911962319095
$ws.Range("T6").Value = 0.1816911962319095

# Row 7
$ws.Range("G7").Value = 0.9943730000000001
$ws.Range("I7").Value = 0.3256940046198011
$ws.Range("J7").Value = 0.325694004619801
$ws.Range("M7").Value = 36.057927
$ws.Range("N7").Value = 108.173781
$ws.Range("O7").Value = 0.1471761890142177
$ws.Range("P7").Value = 0.1471761890142177
$ws.Range("Q7").Value = 35.855029044771
$ws.Range("R7").Value = 322.695261402939
$ws.Range("S7").Value = 0.04793440238472133
$ws.Range("T7").Value = 0.04793440238472133

# Row 8
$ws.Range("G8").Value = 2.018446666666666
$ws.Range("H8").Value = 6.055339999999999
$ws.Range("I8").Value = 0.6611160781499047
$ws.Range("J8").Value = 0.6611160781499047
$ws.Range("M8").Value = 72.266001
$ws.Range("N8").Value = 216.798003
$ws.Range("O8").Value = 0.2949652269937106
$ws.Range("P8").Value = 0.2949652269937106
$ws.Range("Q8").Value = 145.86506883178
$ws.Range("R8").Value = 1312.78561948602
$ws.Range("S8").Value = 0.1950062540606783
$ws.Range("T8").Value = 0.1950062540606783

# Row 9
$ws.Range("G9").Value = 2.018446666666666
$ws.Range("H9").Value = 6.055339999999999
$ws.Range("I9").Value = 0.6611160781499047
$ws.Range("J9").Value = 0.6611160781499047
$ws.Range("N9").Value = 410.023338
$ws.Range("O9").Value = 0.5578585839920717
$ws.Range("P9").Value = 0.5578585839920718
$ws.Range("Q9").Value = 275.8700799472132
$ws.Range("R9").Value = 2482.830719524919
$ws.Range("S9").Value = 0.3688092792110976
$ws.Range("T9").Value = 0.3688092792110977

# Row 10
$ws.Range("G10").Value = 2.018446666666666
$ws.Range("H10").Value = 6.055339999999999
$ws.Range("I10").Value = 0.6611160781499047
$ws.Range("J10").Value = 0.6611160781499047
$ws.Range("M10").Value = 36.057927
$ws.Range("N10").Value = 108.173781
$ws.Range("O10").Value = 0.1471761890142177
$ws.Range("P10").Value = 0.1471761890142177
$ws.Range("Q10").Value = 72.78100256005999
$ws.Range("R10").Value = 655.0290230405399
$ws.Range("S10").Value = 0.0973005448781287
$ws.Range("T10").Value = 0.0973005448781287
